# HW2 req and design
#
# 1) Bump the cached "datetimeFigureOut" date-placeholder text from 10/7/24 to
#    10/8/24 across the slide master and every slide layout (ppPlaceholderDate = 16).
# 2) Rename the "Migration Management" label to "Migration Tracking" and shrink
#    its textbox width on slide 1 and slide 2. The textbox lives inside a group,
#    but Shape.Width/Left/Top/Height are always reported/set in absolute,
#    slide-level coordinates even for shapes reached through GroupItems.

$p = $ppt.ActivePresentation

$oldDate = "10/7/24"
$newDate = "10/8/24"
$ppPlaceholderDate = 16

# --- Slide master date placeholder -----------------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $sh = $masterShapes.Item($j)
    $isDate = $false
    try {
        if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $isDate = $true
        }
    } catch {
    }
    if ($isDate -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout's date placeholder ---------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $layoutShapes = $layout.Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $sh = $layoutShapes.Item($j)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Slide 1 & Slide 2: "Migration Management" textbox (inside a group) ----
for ($slideIdx = 1; $slideIdx -le 2; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $topShape = $slide.Shapes.Item($i)
        if ($topShape.Type -eq 6) {
            for ($k = 1; $k -le $topShape.GroupItems.Count; $k++) {
                $inner = $topShape.GroupItems.Item($k)
                if ($inner.HasTextFrame -and $inner.TextFrame.TextRange.Text -eq "Migration Management") {
                    $inner.TextFrame.TextRange.Text = "Migration Tracking"
                    $inner.Width = 155.615
                }
            }
        }
    }
}

Write-Host "edit complete"
